# Sync attendance_reports, modules_schedules, and assets from main repo - 2026-01-03 18:21:14
#
# The "Recorded By" column (G) on the "Session Analysis Results" sheet lists
# the users who recorded/edited attendance for a session, separated by ", ".
# This reorders the entries in two specific recurring combinations:
#   "System, dnasr281@gmail.com"                 -> "dnasr281@gmail.com, System"
#   "system, backup@backdoor.com, System"        -> "backup@backdoor.com, system, System"
# All other combinations are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($val -eq "System, dnasr281@gmail.com") {
        $cell.Value2 = "dnasr281@gmail.com, System"
    }
    elseif ($val -eq "system, backup@backdoor.com, System") {
        $cell.Value2 = "backup@backdoor.com, system, System"
    }
}
